# Update "想去人数" (people interested) counts in column F across sheets,
# reflecting a refreshed data pull (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4620
$ws1.Range("F6").Value  = 1792
$ws1.Range("F8").Value  = 716
$ws1.Range("F12").Value = 1125
$ws1.Range("F13").Value = 1567
$ws1.Range("F19").Value = 151
$ws1.Range("F21").Value = 1183
$ws1.Range("F25").Value = 1528
$ws1.Range("F30").Value = 66

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 4157
$ws2.Range("F12").Value = 33

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 55

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value  = 4620
$ws4.Range("F11").Value = 1792
$ws4.Range("F12").Value = 716
$ws4.Range("F17").Value = 1125
$ws4.Range("F18").Value = 1567
$ws4.Range("F19").Value = 33
$ws4.Range("F26").Value = 151
$ws4.Range("F33").Value = 1183
$ws4.Range("F43").Value = 1528
$ws4.Range("F50").Value = 66
